# Update currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H, I, J, K, L, M, N) on affected leve rows across multiple sheets,
# reflecting refreshed market-board prices from the scheduled data-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1249.5
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 1499
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 4497
$ws.Range("M46").Value = -2881
$ws.Range("N46").Value = -4735

$ws.Range("H60").Value = 1249.5
$ws.Range("I60").Value = 1000
$ws.Range("J60").Value = 1499
$ws.Range("K60").Value = 3000
$ws.Range("L60").Value = 4497
$ws.Range("M60").Value = -2516
$ws.Range("N60").Value = -5465

$ws.Range("H74").Value = 6873.6313
$ws.Range("I74").Value = 3958.8
$ws.Range("J74").Value = 7914.643
$ws.Range("K74").Value = 3958.8
$ws.Range("L74").Value = 7914.643
$ws.Range("M74").Value = -3022.8
$ws.Range("N74").Value = -9786.643

$ws.Range("H77").Value = 6873.6313
$ws.Range("I77").Value = 3958.8
$ws.Range("J77").Value = 7914.643
$ws.Range("K77").Value = 19794
$ws.Range("L77").Value = 39573.215
$ws.Range("M77").Value = -15114
$ws.Range("N77").Value = -48933.215

$ws.Range("H100").Value = 3469.5
$ws.Range("I100").Value = 3469.5
$ws.Range("K100").Value = 3469.5
$ws.Range("M100").Value = -2928.5

$ws.Range("H116").Value = 4983.6
$ws.Range("J116").Value = 4982.3335
$ws.Range("L116").Value = 4982.3335
$ws.Range("N116").Value = -11866.3335

$ws.Range("H135").Value = 1144.5172
$ws.Range("I135").Value = 913.13635
$ws.Range("K135").Value = 8218.227150000001
$ws.Range("M135").Value = -5683.227150000001

$ws.Range("H136").Value = 199772.75
$ws.Range("J136").Value = 199772.75
$ws.Range("L136").Value = 199772.75
$ws.Range("N136").Value = -209972.75

$ws.Range("H137").Value = 4280.758
$ws.Range("J137").Value = 4728.905
$ws.Range("L137").Value = 14186.715
$ws.Range("N137").Value = -19286.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 15351.143
$ws.Range("J26").Value = 55000
$ws.Range("L26").Value = 55000
$ws.Range("N26").Value = -55584

$ws.Range("H28").Value = 50000
$ws.Range("J28").Value = 50000
$ws.Range("L28").Value = 50000
$ws.Range("N28").Value = -50588

$ws.Range("H134").Value = 2568.8164
$ws.Range("I134").Value = 888.2222
$ws.Range("K134").Value = 2664.6666
$ws.Range("M134").Value = -129.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 750.5714
$ws.Range("J10").Value = 899
$ws.Range("L10").Value = 899
$ws.Range("N10").Value = -1177

$ws.Range("H31").Value = 2535.074
$ws.Range("I31").Value = 1936.2693
$ws.Range("J31").Value = 3091.1072
$ws.Range("K31").Value = 1936.2693
$ws.Range("L31").Value = 3091.1072
$ws.Range("M31").Value = -1641.2693
$ws.Range("N31").Value = -3681.1072

$ws.Range("H34").Value = 2535.074
$ws.Range("I34").Value = 1936.2693
$ws.Range("J34").Value = 3091.1072
$ws.Range("K34").Value = 1936.2693
$ws.Range("L34").Value = 3091.1072
$ws.Range("M34").Value = -1734.2693
$ws.Range("N34").Value = -3495.1072

$ws.Range("H134").Value = 19345.887
$ws.Range("I134").Value = 23772.268
$ws.Range("K134").Value = 71316.804
$ws.Range("M134").Value = -68781.804

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1274.125
$ws.Range("I68").Value = 598.5
$ws.Range("J68").Value = 1949.75
$ws.Range("K68").Value = 1795.5
$ws.Range("L68").Value = 5849.25
$ws.Range("M68").Value = -984.5
$ws.Range("N68").Value = -7471.25

$ws.Range("H71").Value = 1274.125
$ws.Range("I71").Value = 598.5
$ws.Range("J71").Value = 1949.75
$ws.Range("K71").Value = 5386.5
$ws.Range("L71").Value = 17547.75
$ws.Range("M71").Value = -1330.5
$ws.Range("N71").Value = -25659.75

$ws.Range("H98").Value = 621.5
$ws.Range("I98").Value = 222
$ws.Range("K98").Value = 666
$ws.Range("M98").Value = 832

$ws.Range("H109").Value = 1912.6
$ws.Range("I109").Value = 1140.75
$ws.Range("J109").Value = 5000
$ws.Range("K109").Value = 3422.25
$ws.Range("L109").Value = 15000
$ws.Range("M109").Value = -2382.25
$ws.Range("N109").Value = -17080

$ws.Range("H136").Value = 2404
$ws.Range("I136").Value = 2404
$ws.Range("K136").Value = 7212
$ws.Range("M136").Value = -2112

$ws.Range("H139").Value = 2774.5
$ws.Range("I139").Value = 2750
$ws.Range("J139").Value = 2799
$ws.Range("K139").Value = 8250
$ws.Range("L139").Value = 8397
$ws.Range("M139").Value = -3110
$ws.Range("N139").Value = -18677

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 79852.586
$ws.Range("I122").Value = 99129.58
$ws.Range("K122").Value = 297388.74
$ws.Range("M122").Value = -294938.74

$ws.Range("H132").Value = 3196.44
$ws.Range("J132").Value = 5062.125
$ws.Range("L132").Value = 15186.375
$ws.Range("N132").Value = -20246.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 745
$ws.Range("I3").Value = 775
$ws.Range("J3").Value = 715
$ws.Range("K3").Value = 775
$ws.Range("L3").Value = 715
$ws.Range("M3").Value = -663
$ws.Range("N3").Value = -939

$ws.Range("H15").Value = 745
$ws.Range("I15").Value = 775
$ws.Range("J15").Value = 715
$ws.Range("K15").Value = 775
$ws.Range("L15").Value = 715
$ws.Range("M15").Value = -605
$ws.Range("N15").Value = -1055

$ws.Range("H20").Value = 7750
$ws.Range("I20").Value = 7166.6665
$ws.Range("J20").Value = 7866.6665
$ws.Range("K20").Value = 7166.6665
$ws.Range("L20").Value = 7866.6665
$ws.Range("M20").Value = -6940.6665
$ws.Range("N20").Value = -8318.666499999999

$ws.Range("H22").Value = 1170.4166
$ws.Range("J22").Value = 1899.3334
$ws.Range("L22").Value = 1899.3334
$ws.Range("N22").Value = -2489.3334

$ws.Range("H27").Value = 1170.4166
$ws.Range("J27").Value = 1899.3334
$ws.Range("L27").Value = 1899.3334
$ws.Range("N27").Value = -2113.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 964.0540999999999
$ws.Range("I113").Value = 763
$ws.Range("K113").Value = 2289
$ws.Range("M113").Value = -119

$ws.Range("H132").Value = 171412.08
$ws.Range("I132").Value = 1818.2727
$ws.Range("J132").Value = 2503327
$ws.Range("K132").Value = 5454.8181
$ws.Range("L132").Value = 7509981
$ws.Range("M132").Value = -2924.8181
$ws.Range("N132").Value = -7515041
